$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grids")

$lastRow = $ws.UsedRange.Rows.Count
$firstDataRow = 5
$col = 3

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "" -and $val.ToString().StartsWith("e_") -eq $false) {
        $cell.Value = "e_" + $val
    }
}

$ws.Columns.Item($col).ColumnWidth = 14.5

"Renamed grid node codes in column C with e_ prefix"
